$wb = $excel.ActiveWorkbook
$wsSound = $wb.Worksheets.Item("SOUND")
$wsVoice = $wb.Worksheets.Item("VOICE")

# --- SOUND sheet: key_bundle column (D) "sound" -> "sounds" for data rows 5-11 ---
foreach ($r in 5..11) {
    $wsSound.Cells.Item($r, 4).Value = "sounds"
}

# --- VOICE sheet: key_bundle column (B) "greet" -> "voices" for data rows 5-6 ---
$wsVoice.Range("B5").Value = "voices"
$wsVoice.Range("B6").Value = "voices"

# --- Auto-fit column B on VOICE sheet (width became 11.25 due to "voices" content) ---
$wsVoice.Columns.Item(2).AutoFit()

# --- Select C6 on SOUND sheet, then activate VOICE sheet and select it last ---
$wsSound.Activate()
$wsSound.Range("C6").Select()

$wsVoice.Activate()
